$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate  = "10-11-2025"
$newPrice = "The price of gold in India today is ₹12,322 per gram for 24 karat gold, ₹11,295 per gram for 22 karat gold and ₹9,242 per gram for 18 karat gold (also called 999 gold)."

# Bring formatting for the new row in line with the existing data rows
# (thin border via style col A / border+wrap style col B, same as row 54)
# before the cell values are entered.
$ws.Range("A54:B54").Copy()
$ws.Range("A55:B55").PasteSpecial(-4122)

# Enter the date as literal text (leading apostrophe forces text so Excel
# doesn't reinterpret "10-11-2025" as a date serial), then re-apply the
# plain border formatting from row 54 so the cell collapses back onto the
# same style as the rest of the date column instead of keeping a
# quote-prefixed variant.
$ws.Range("A55").Value = "'" + $newDate
$ws.Range("A54").Copy()
$ws.Range("A55").PasteSpecial(-4122)

$ws.Range("B55").Value = $newPrice

$excel.CutCopyMode = 0
